# PowerPoint COM-interop script
# Rewrites the tokistorage-client-deck.pptx text content per the commit:
# 'client deck PPTX/PDF を client-proposal.html と完全同期' — reframing all 10
# slides around '未来の自分からの言葉・ブレない軸' instead of '100年後に何を残すか'.
#
# NOTE: we always replace text via TextRange.Characters(1, <old length>).Text = ...
# rather than TextRange.Text = ... directly. Both exist on the object model, but
# assigning the whole-range .Text property flattens an embedded line-break into a
# brand-new <a:p> paragraph (losing the original single-run/embedded-\n shape of the
# run). Re-assigning through .Characters(...) targets the run in place and keeps a
# literal line-break character inside the single <a:t> run, matching the source deck's
# existing convention (and the target diff) exactly.

$p = $ppt.ActivePresentation

# ---- Slide 1 ----
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(5).TextFrame.TextRange  # TextBox 5
$tr.Characters(1, $tr.Length).Text = "未来の自分は、`n今のあなたに何を言いますか。"
$tr = $s.Shapes.Item(6).TextFrame.TextRange  # TextBox 6
$tr.Characters(1, $tr.Length).Text = "AIがあらゆる問いに答える時代。`nしかしこの問いだけは、あなた自身にしか答えられない。`nその言葉を三層に永続化することで、ブレない軸が生まれます。"

# ---- Slide 2 ----
$s = $p.Slides.Item(2)
$tr = $s.Shapes.Item(2).TextFrame.TextRange  # TextBox 2
$tr.Characters(1, $tr.Length).Text = "AIが答えられない、ただひとつの問い──未来の自分の声"
$tr = $s.Shapes.Item(6).TextFrame.TextRange  # TextBox 6
$tr.Characters(1, $tr.Length).Text = "AIは「最適解」を出す。しかし「あなたの軸」は決められない。"
$tr = $s.Shapes.Item(7).TextFrame.TextRange  # TextBox 7
$tr.Characters(1, $tr.Length).Text = "生成AIは調査・分析・戦略を瞬時に出力します。しかし「未来の自分は今の自分に何を言うか」——この問いへの答えは、あなたの人生の文脈からしか生まれません。最適化ではなく、軸の選択です。"
$tr = $s.Shapes.Item(10).TextFrame.TextRange  # TextBox 10
$tr.Characters(1, $tr.Length).Text = "目標は消えるから、ブレる。"
$tr = $s.Shapes.Item(11).TextFrame.TextRange  # TextBox 11
$tr.Characters(1, $tr.Length).Text = "手帳に書いた抱負は年末に忘れている。デジタルのメモは通知の洪水に埋もれる。人が軸を保てないのは、意志が弱いからではない。言葉が消えるからです。消えない形で刻まれた言葉だけが、軸として機能し続けます。"
$tr = $s.Shapes.Item(14).TextFrame.TextRange  # TextBox 14
$tr.Characters(1, $tr.Length).Text = "「未来の自分」は、脳にとって他人である。"
$tr = $s.Shapes.Item(15).TextFrame.TextRange  # TextBox 15
$tr.Characters(1, $tr.Length).Text = "心理学者ハーシュフィールドのfMRI研究によれば、「未来の自分」は脳の活動パターンにおいて「他人」と同じように処理されます。しかし「千年残る」という重みが、未来の自分の臨場感を引き上げる。永続記録の不可逆性が、未来の自分を「他人」から「自分」に引き戻します。"

# ---- Slide 3 ----
$s = $p.Slides.Item(3)
$tr = $s.Shapes.Item(2).TextFrame.TextRange  # TextBox 2
$tr.Characters(1, $tr.Length).Text = "あなたが買うのは「プロダクト」ではありません。未来の自分と出会うための対話です。"
$tr = $s.Shapes.Item(6).TextFrame.TextRange  # TextBox 6
$tr.Characters(1, $tr.Length).Text = "未来の自分の声を聴き、`nブレない軸を永続化する。"
$tr = $s.Shapes.Item(9).TextFrame.TextRange  # TextBox 9
$tr.Characters(1, $tr.Length).Text = "未来の自分との対話"
$tr = $s.Shapes.Item(10).TextFrame.TextRange  # TextBox 10
$tr.Characters(1, $tr.Length).Text = "「未来の自分は、今の自分に何を言うか」——この問いに、9領域の思想フレームワークで向き合います。浮かんだ言葉は、あなただけのものです。"
$tr = $s.Shapes.Item(13).TextFrame.TextRange  # TextBox 13
$tr.Characters(1, $tr.Length).Text = "軸の言語化"
$tr = $s.Shapes.Item(14).TextFrame.TextRange  # TextBox 14
$tr.Characters(1, $tr.Length).Text = "対話から浮かんだ言葉を、あなたの人生の物語・価値観とともに構造化。「ブレない軸」として言語化し、千年先まで届く形に設計します。"
$tr = $s.Shapes.Item(17).TextFrame.TextRange  # TextBox 17
$tr.Characters(1, $tr.Length).Text = "ブレない軸の永続化"
$tr = $s.Shapes.Item(18).TextFrame.TextRange  # TextBox 18
$tr.Characters(1, $tr.Length).Text = "言語化された軸を、消えない形で刻む。千年残るという重みが不可逆の約束になる。消える言葉では軸にならない。トキストレージだからこそ成立する、自己変容の装置。"
$tr = $s.Shapes.Item(21).TextFrame.TextRange  # TextBox 21
$tr.Characters(1, $tr.Length).Text = "三層分散保管への格納"
$tr = $s.Shapes.Item(22).TextFrame.TextRange  # TextBox 22
$tr.Characters(1, $tr.Length).Text = "あなたの軸と存在証明を三層に分散保管。物理層（石英ガラス）、国家層（国立国会図書館）、民間層（GitHub）──単一障害点のない構造で、千年先まで届けます。"

# ---- Slide 4 ----
$s = $p.Slides.Item(4)
$tr = $s.Shapes.Item(2).TextFrame.TextRange  # TextBox 2
$tr.Characters(1, $tr.Length).Text = "対話から始まり、三層の分散保管で千年に届く"
$tr = $s.Shapes.Item(8).TextFrame.TextRange  # TextBox 8
$tr.Characters(1, $tr.Length).Text = "「未来の自分は、今の自分に何を言うか」を軸に、あなたの人生の物語を聴かせていただきます。この対話自体が、多くの方にとって深い内省の機会になります。"
$tr = $s.Shapes.Item(12).TextFrame.TextRange  # TextBox 12
$tr.Characters(1, $tr.Length).Text = "ブレない軸の言語化"
$tr = $s.Shapes.Item(13).TextFrame.TextRange  # TextBox 13
$tr.Characters(1, $tr.Length).Text = "未来の自分を想像し、その視座から今の自分への言葉を紡ぎます。存在証明の構造設計と同時に、あなたの「ブレない軸」を言葉にするプロセスです。"
$tr = $s.Shapes.Item(17).TextFrame.TextRange  # TextBox 17
$tr.Characters(1, $tr.Length).Text = "コンテンツの制作"
$tr = $s.Shapes.Item(18).TextFrame.TextRange  # TextBox 18
$tr.Characters(1, $tr.Length).Text = "設計に基づいて、あなたの存在証明を制作します。必要に応じて、プロフェッショナルな撮影・収録・編集をコーディネート。あなたの物語と軸を最高の形に仕上げます。"
$tr = $s.Shapes.Item(22).TextFrame.TextRange  # TextBox 22
$tr.Characters(1, $tr.Length).Text = "三層分散保管・納品"
$tr = $s.Shapes.Item(23).TextFrame.TextRange  # TextBox 23
$tr.Characters(1, $tr.Length).Text = "存在証明とブレない軸を三層に格納。石英ガラスへの刻印（物理層）、国立国会図書館への納本（国家層）、GitHubへの分散保管（民間層）。QRコードを読み取るたびに、未来の自分があなたに語りかけます。"

# ---- Slide 5 ----
$s = $p.Slides.Item(5)
$tr = $s.Shapes.Item(2).TextFrame.TextRange  # TextBox 2
$tr.Characters(1, $tr.Length).Text = "ブレない軸を必要としている、すべての方へ"
$tr = $s.Shapes.Item(3).TextFrame.TextRange  # TextBox 3
$tr.Characters(1, $tr.Length).Text = "WHO THIS IS TOR"
$tr = $s.Shapes.Item(8).TextFrame.TextRange  # TextBox 8
$tr.Characters(1, $tr.Length).Text = "事業の浮き沈みに揺れない、自分自身の理念を軸として永続化したい"
$tr = $s.Shapes.Item(11).TextFrame.TextRange  # TextBox 11
$tr.Characters(1, $tr.Length).Text = "T"
$tr = $s.Shapes.Item(12).TextFrame.TextRange  # TextBox 12
$tr.Characters(1, $tr.Length).Text = "人生の転機にいる方"
$tr = $s.Shapes.Item(13).TextFrame.TextRange  # TextBox 13
$tr.Characters(1, $tr.Length).Text = "転職、独立、退職——次の一歩を踏み出す前に、ブレない軸を言語化したい"
$tr = $s.Shapes.Item(17).TextFrame.TextRange  # TextBox 17
$tr.Characters(1, $tr.Length).Text = "家族に軸を残したい方"
$tr = $s.Shapes.Item(18).TextFrame.TextRange  # TextBox 18
$tr.Characters(1, $tr.Length).Text = "子や孫に伝えたいのは財産ではなく、生き方の軸。消えない形で届けたい"
$tr = $s.Shapes.Item(22).TextFrame.TextRange  # TextBox 22
$tr.Characters(1, $tr.Length).Text = "アーティスト・クリエイター"
$tr = $s.Shapes.Item(23).TextFrame.TextRange  # TextBox 23
$tr.Characters(1, $tr.Length).Text = "創作の原点を見失わないために、自分の軸をプラットフォームに依存せず刻みたい"
$tr = $s.Shapes.Item(27).TextFrame.TextRange  # TextBox 27
$tr.Characters(1, $tr.Length).Text = "宗教者・教育者"
$tr = $s.Shapes.Item(28).TextFrame.TextRange  # TextBox 28
$tr.Characters(1, $tr.Length).Text = "教えの本質を、自分がいなくなっても消えない形で次の世代に手渡したい"
$tr = $s.Shapes.Item(32).TextFrame.TextRange  # TextBox 32
$tr.Characters(1, $tr.Length).Text = "地域・コミュニティ"
$tr = $s.Shapes.Item(33).TextFrame.TextRange  # TextBox 33
$tr.Characters(1, $tr.Length).Text = "まちの精神、災害の教訓、創設の志──組織の軸を人の寿命から解放したい"

# ---- Slide 6 ----
$s = $p.Slides.Item(6)
$tr = $s.Shapes.Item(2).TextFrame.TextRange  # TextBox 2
$tr.Characters(1, $tr.Length).Text = "対話から生まれ、三層に分散保管される具体的な成果物"
$tr = $s.Shapes.Item(12).TextFrame.TextRange  # TextBox 12
$tr.Characters(1, $tr.Length).Text = "未来の自分からの言葉──ブレない軸の永続化"
$tr = $s.Shapes.Item(13).TextFrame.TextRange  # TextBox 13
$tr.Characters(1, $tr.Length).Text = "「千年残る」という重みが、未来の自分の臨場感を引き上げる。その未来の自分から今のあなたへ投げかけられた言葉を、三層に永続化する。消えないからこそ、ブレない軸になる。"
$tr = $s.Shapes.Item(18).TextFrame.TextRange  # TextBox 18
$tr.Characters(1, $tr.Length).Text = "あなたの存在証明を、70以上の思想エッセイの文脈に位置づけたレポート。「なぜ残すのか」の知的な裏付けを提供します。"

# ---- Slide 7 ----
$s = $p.Slides.Item(7)
$tr = $s.Shapes.Item(2).TextFrame.TextRange  # TextBox 2
$tr.Characters(1, $tr.Length).Text = "未来の自分からの言葉を、どこまで深く刻むか"
$tr = $s.Shapes.Item(8).TextFrame.TextRange  # TextBox 8
$tr.Characters(1, $tr.Length).Text = "「未来の自分は何を言うか」──その最初の一言を声で刻み、三層に格納する体験。A4ラミネート＋国立国会図書館納本＋GitHub。まだ軸が見えなくても、ここから始まります。"
$tr = $s.Shapes.Item(13).TextFrame.TextRange  # TextBox 13
$tr.Characters(1, $tr.Length).Text = "対話から生まれたブレない軸を、石英ガラスに刻む。千年残る不可逆の約束が、日常の中で軸として機能し続ける。年あたり50円の自己変容装置。"
$tr = $s.Shapes.Item(18).TextFrame.TextRange  # TextBox 18
$tr.Characters(1, $tr.Length).Text = "三世代先の子孫から今の自分に投げかける言葉を刻む。その言葉を受け取った瞬間、未来世代への責任が自覚になる。お墓や仏壇は風化するが、軸は消えない。世代を超えて対話が続く。"
$tr = $s.Shapes.Item(23).TextFrame.TextRange  # TextBox 23
$tr.Characters(1, $tr.Length).Text = "声そのものを千年先に届ける。音声復元技術を含む完全オーダーメイド。未来の子孫がQRコードを読み取ったとき、あなたの肉声が語りかける。"

# ---- Slide 9 ----
$s = $p.Slides.Item(9)
$tr = $s.Shapes.Item(7).TextFrame.TextRange  # TextBox 7
$tr.Characters(1, $tr.Length).Text = "「未来の自分は今の自分に何を言うか」を一緒に探る、最初の90分"
$tr = $s.Shapes.Item(10).TextFrame.TextRange  # TextBox 10
$tr.Characters(1, $tr.Length).Text = "ブレない軸の言語化"
$tr = $s.Shapes.Item(11).TextFrame.TextRange  # TextBox 11
$tr.Characters(1, $tr.Length).Text = "未来の自分の視座から、今の自分への言葉を紡ぎます"
$tr = $s.Shapes.Item(14).TextFrame.TextRange  # TextBox 14
$tr.Characters(1, $tr.Length).Text = "制作・三層分散保管"
$tr = $s.Shapes.Item(15).TextFrame.TextRange  # TextBox 15
$tr.Characters(1, $tr.Length).Text = "存在証明とブレない軸を、物理・国家・民間の三層に格納します"
$tr = $s.Shapes.Item(18).TextFrame.TextRange  # TextBox 18
$tr.Characters(1, $tr.Length).Text = "納品"
$tr = $s.Shapes.Item(19).TextFrame.TextRange  # TextBox 19
$tr.Characters(1, $tr.Length).Text = "QRコードを読み取るたびに、未来の自分が語りかけます"
